$d = $word.ActiveDocument
$nbsp = [char]0xA0

# --- Remove the two HYPERLINK fields for the CryptoBridge / HitBTC exchange
# mentions (they are Word fields: fldChar begin/instrText/separate/<result>/end).
# Field.Delete() removes the whole field construct (begin..end) as a unit,
# matching the removal of those runs in the diff.
$fieldsToDelete = @()
foreach ($f in $d.Fields) {
  if ($f.Code.Text -like "*crypto-bridge.org*" -or $f.Code.Text -like "*hitbtc.com*") {
    $fieldsToDelete += $f
  }
}
foreach ($f in $fieldsToDelete) {
  $f.Delete()
}

# --- Trim "... exchanges such as<nbsp>" down to "... exchanges" (the bullet's
# lead-in sentence now ends right before the spot the removed fields used to
# occupy).
$d.Content.Find.Execute("exchanges such as" + $nbsp, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "exchanges", 2) | Out-Null

# --- The ", " run that used to sit between the two now-deleted hyperlink
# fields is orphaned; delete just that little run's range (not a text
# replace) so the bold "." run right after it keeps its own identity/run
# formatting instead of getting merged into a neighboring run.
$sepRange = $d.Content
if ($sepRange.Find.Execute("," + $nbsp, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)) {
  $sepRange.Delete()
}

# --- Renumber the "smarthosting" bookmark's underlying w:id back to 0 (Word
# recomputes bookmark ids on save; deleting + re-adding the bookmark in place
# forces that recompute deterministically).
$bm = $d.Bookmarks.Item("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange) | Out-Null
